$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell forcing TEXT storage (no numeric
# auto-conversion, no style/number-format drift on the destination cell).
# Uses a scratch cell far outside the used range, formatted as Text,
# then Copy / PasteSpecial(values-only) into the target so the target
# keeps its original (default) style index.
function Set-TextValue($cellRef, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

Set-TextValue "D2" "56.502.39"
$ws.Range("E2").Value = "  -0.93%  "

Set-TextValue "D3" "2.330.84"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  +0.14%  "

Set-TextValue "D5" "513.02"
$ws.Range("E5").Value = "  -1.27%  "

Set-TextValue "D6" "132.15"
$ws.Range("E6").Value = "  -2.02%  "

$ws.Range("E7").Value = "  -0.05%  "

Set-TextValue "D8" "0.534"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("E9").Value = "  -3.25%  "

$ws.Range("E10").Value = "  -0.51%  "

$ws.Range("E11").Value = "  +1.82%  "

$ws.Range("E12").Value = "  -0.64%  "

Set-TextValue "D13" "2.746.75"
$ws.Range("E13").Value = "  -0.56%  "

Set-TextValue "D14" "23.56"
$ws.Range("E14").Value = "  -0.66%  "

Set-TextValue "D15" "56.487.25"
$ws.Range("E15").Value = "  -0.69%  "

Set-TextValue "D16" "0.0000133"
$ws.Range("E16").Value = "  -1.44%  "

Set-TextValue "D17" "2.326.50"
$ws.Range("E17").Value = "  -1.55%  "

$ws.Range("E18").Value = "  -0.27%  "

Set-TextValue "D19" "324.53"
$ws.Range("E19").Value = "  +0.45%  "

Set-TextValue "D20" "4.13"
$ws.Range("E20").Value = "  -2.65%  "

Set-TextValue "D21" "6.66"
$ws.Range("E21").Value = "  +1.81%  "

$ws.Range("E22").Value = "  -0.14%  "

Set-TextValue "D23" "61.85"
$ws.Range("E23").Value = "  +1.41%  "

Set-TextValue "D24" "8.68"
$ws.Range("E24").Value = "  +11.58%  "

$ws.Range("E25").Value = "  +1.18%  "

$ws.Range("E26").Value = "  -7.15%  "

$ws.Range("E27").Value = "  +5.07%  "

Set-TextValue "D28" "167.70"
$ws.Range("E28").Value = "  -1.48%  "

$ws.Range("E29").Value = "  -1.21%  "

Set-TextValue "D30" "0.0₃0721"
$ws.Range("E30").Value = "  -3.38%  "

Set-TextValue "D31" "6.11"
$ws.Range("E31").Value = "  -1.37%  "

Set-TextValue "D32" "18.33"
$ws.Range("E32").Value = "  +0.30%  "

$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("E36").Value = "  -1.92%  "

Set-TextValue "D37" "0.889"
$ws.Range("E37").Value = "  -4.49%  "

$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D38" "153.42"
$ws.Range("E38").Value = "  +11.77%  "

Set-TextValue "D39" "1.56"
$ws.Range("E39").Value = "  +1.08%  "

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D40" "38.45"
$ws.Range("E40").Value = "  +1.67%  "

Set-TextValue "D41" "0.376"
$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("E42").Value = "  -0.74%  "

Set-TextValue "D43" "279.01"
$ws.Range("E43").Value = "  +0.38%  "

Set-TextValue "D44" "5.04"
$ws.Range("E44").Value = "  -1.94%  "

Set-TextValue "D45" "0.0927"
$ws.Range("E45").Value = "  -0.88%  "

Set-TextValue "D46" "0.0496"
$ws.Range("E46").Value = "  -1.73%  "

Set-TextValue "D47" "0.559"
$ws.Range("E47").Value = "  -0.77%  "

Set-TextValue "D48" "18.12"
$ws.Range("E48").Value = "  +5.15%  "

$ws.Range("B49").Value = "Polygon"
$ws.Range("C49").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D49" "0.382"
$ws.Range("E49").Value = "  +1.05%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D50" "0.0214"
$ws.Range("E50").Value = "  -1.88%  "

Set-TextValue "D51" "17.13"
$ws.Range("E51").Value = "  +2.05%  "

$excel.CutCopyMode = 0
